$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 903.1905
$ws.Range("I2").Value = 221.21428
$ws.Range("K2").Value = 221.21428
$ws.Range("M2").Value = -108.21428

$ws.Range("H9").Value = 409.2857
$ws.Range("I9").Value = 430.41177
$ws.Range("J9").Value = 319.5
$ws.Range("K9").Value = 430.41177
$ws.Range("L9").Value = 319.5
$ws.Range("M9").Value = -261.41177
$ws.Range("N9").Value = -657.5

$ws.Range("H101").Value = 16667157
$ws.Range("I101").Value = 25000436
$ws.Range("K101").Value = 75001308
$ws.Range("M101").Value = -74999686

$ws.Range("H116").Value = 4260.7
$ws.Range("I116").Value = 1800.4
$ws.Range("J116").Value = 6721
$ws.Range("K116").Value = 1800.4
$ws.Range("L116").Value = 6721
$ws.Range("M116").Value = 1641.6
$ws.Range("N116").Value = -13605

$ws.Range("H132").Value = 3730.3
$ws.Range("I132").Value = 716.6667
$ws.Range("K132").Value = 2150.0001
$ws.Range("M132").Value = 379.9998999999998

$ws.Range("H138").Value = 2638.5881
$ws.Range("J138").Value = 2097.3333
$ws.Range("L138").Value = 6291.999899999999
$ws.Range("N138").Value = -16571.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1949.3334
$ws.Range("I21").Value = 1949.3334
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1949.3334
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -1575.3334
$ws.Range("N21").ClearContents()

$ws.Range("H24").Value = 38339.8
$ws.Range("J24").Value = 38339.8
$ws.Range("L24").Value = 38339.8
$ws.Range("N24").Value = -39087.8

$ws.Range("H32").Value = 3200.5938
$ws.Range("I32").Value = 2747.4666
$ws.Range("J32").Value = 9997.5
$ws.Range("K32").Value = 2747.4666
$ws.Range("L32").Value = 9997.5
$ws.Range("M32").Value = -2460.4666
$ws.Range("N32").Value = -10571.5

$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H61").Value = 4729.45
$ws.Range("I61").Value = 3709
$ws.Range("K61").Value = 3709
$ws.Range("M61").Value = -3497

$ws.Range("H82").Value = 22794.75
$ws.Range("J82").Value = 22794.75
$ws.Range("L82").Value = 22794.75
$ws.Range("N82").Value = -23516.75

$ws.Range("H85").Value = 22794.75
$ws.Range("J85").Value = 22794.75
$ws.Range("L85").Value = 22794.75
$ws.Range("N85").Value = -25290.75

$ws.Range("H100").Value = 38339.8
$ws.Range("J100").Value = 38339.8
$ws.Range("L100").Value = 38339.8
$ws.Range("N100").Value = -40503.8

$ws.Range("H122").Value = 3617.8
$ws.Range("I122").Value = 3430
$ws.Range("J122").Value = 3899.5
$ws.Range("K122").Value = 10290
$ws.Range("L122").Value = 11698.5
$ws.Range("M122").Value = -7840
$ws.Range("N122").Value = -16598.5

$ws.Range("H136").Value = 4729.45
$ws.Range("I136").Value = 3709
$ws.Range("K136").Value = 11127
$ws.Range("M136").Value = -8577

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 66430
$ws.Range("J19").Value = 66430
$ws.Range("L19").Value = 66430
$ws.Range("N19").Value = -66776

$ws.Range("H20").Value = 1087.75
$ws.Range("J20").Value = 950.75
$ws.Range("L20").Value = 950.75
$ws.Range("N20").Value = -1444.75

$ws.Range("H76").Value = 200314
$ws.Range("J76").Value = 200314
$ws.Range("L76").Value = 200314
$ws.Range("N76").Value = -200944

$ws.Range("H79").Value = 200314
$ws.Range("J79").Value = 200314
$ws.Range("L79").Value = 200314
$ws.Range("N79").Value = -202498

$ws.Range("H88").Value = 13860
$ws.Range("I88").Value = 4000
$ws.Range("J88").Value = 16325
$ws.Range("K88").Value = 4000
$ws.Range("L88").Value = 16325
$ws.Range("M88").Value = -3594
$ws.Range("N88").Value = -17137

$ws.Range("H91").Value = 13860
$ws.Range("I91").Value = 4000
$ws.Range("J91").Value = 16325
$ws.Range("K91").Value = 4000
$ws.Range("L91").Value = 16325
$ws.Range("M91").Value = -2596
$ws.Range("N91").Value = -19133

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 230.2
$ws.Range("I7").Value = 70.5
$ws.Range("K7").Value = 70.5
$ws.Range("M7").Value = 42.5

$ws.Range("H31").Value = 2839.5
$ws.Range("I31").Value = 2038
$ws.Range("K31").Value = 2038
$ws.Range("M31").Value = -1743

$ws.Range("H34").Value = 2839.5
$ws.Range("I34").Value = 2038
$ws.Range("K34").Value = 2038
$ws.Range("M34").Value = -1836

$ws.Range("H92").Value = 40500
$ws.Range("J92").Value = 40500
$ws.Range("L92").Value = 40500
$ws.Range("N92").Value = -45492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4999.5
$ws.Range("I81").Value = 4999
$ws.Range("K81").Value = 14997
$ws.Range("M81").Value = -13874

$ws.Range("H84").Value = 4999.5
$ws.Range("I84").Value = 4999
$ws.Range("K84").Value = 44991
$ws.Range("M84").Value = -39375

$ws.Range("H107").Value = 513.25
$ws.Range("J107").Value = 558.3333
$ws.Range("L107").Value = 1674.9999
$ws.Range("N107").Value = -5514.9999

$ws.Range("H131").Value = 627751.4
$ws.Range("J131").Value = 912548.4
$ws.Range("L131").Value = 2737645.2
$ws.Range("N131").Value = -2747725.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 183.5
$ws.Range("J2").Value = 500
$ws.Range("L2").Value = 500
$ws.Range("N2").Value = -726

$ws.Range("H70").Value = 6920.375
$ws.Range("J70").Value = 7477
$ws.Range("L70").Value = 7477
$ws.Range("N70").Value = -8017

$ws.Range("H73").Value = 6920.375
$ws.Range("J73").Value = 7477
$ws.Range("L73").Value = 7477
$ws.Range("N73").Value = -9349

$ws.Range("H95").Value = 27998
$ws.Range("J95").Value = 27998
$ws.Range("L95").Value = 27998
$ws.Range("N95").Value = -33490

$ws.Range("H97").Value = 843
$ws.Range("I97").Value = 843
$ws.Range("K97").Value = 843
$ws.Range("M97").Value = -347

$ws.Range("H98").Value = 50643
$ws.Range("J98").Value = 50643
$ws.Range("L98").Value = 50643
$ws.Range("N98").Value = -56633

$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 2500
$ws.Range("K102").Value = 2500
$ws.Range("M102").Value = -878

$ws.Range("H122").Value = 2099.5
$ws.Range("I122").Value = 2099
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 6297
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -3847
$ws.Range("N122").Value = -11200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 712.5
$ws.Range("I9").Value = 544
$ws.Range("J9").Value = 1049.5
$ws.Range("K9").Value = 544
$ws.Range("L9").Value = 1049.5
$ws.Range("M9").Value = -320
$ws.Range("N9").Value = -1497.5

$ws.Range("H22").Value = 1159.9166
$ws.Range("I22").Value = 1172
$ws.Range("K22").Value = 1172
$ws.Range("M22").Value = -877

$ws.Range("H27").Value = 1159.9166
$ws.Range("I27").Value = 1172
$ws.Range("K27").Value = 1172
$ws.Range("M27").Value = -1065

$ws.Range("H35").Value = 577
$ws.Range("I35").Value = 577
$ws.Range("K35").Value = 577
$ws.Range("M35").Value = -241

$ws.Range("H40").Value = 4226
$ws.Range("I40").Value = 3752
$ws.Range("K40").Value = 3752
$ws.Range("M40").Value = -3616

$ws.Range("H46").Value = 2797.5
$ws.Range("I46").Value = 2095
$ws.Range("J46").Value = 3500
$ws.Range("K46").Value = 2095
$ws.Range("L46").Value = 3500
$ws.Range("M46").Value = -1907
$ws.Range("N46").Value = -3876

$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10450

$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11560

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H122").Value = 7838.3887
$ws.Range("I122").Value = 9410
$ws.Range("J122").Value = 5873.875
$ws.Range("K122").Value = 28230
$ws.Range("L122").Value = 17621.625
$ws.Range("M122").Value = -25780
$ws.Range("N122").Value = -22521.625

$ws.Range("H132").Value = 1793
$ws.Range("I132").Value = 1793
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5379
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2849
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 23596.75
$ws.Range("J98").Value = 23596.75
$ws.Range("L98").Value = 23596.75
$ws.Range("N98").Value = -29586.75

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
